# EFC_QTR_FIN.xlsx - "Doing Updates for Financials"
# A new reporting quarter (and the prior quarter) was inserted as the two new
# leading data columns (D & E) in front of the existing quarterly columns on
# the "EFC" sheet; all previously-existing quarter columns (old D:K) shift
# right to F:M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert two new blank columns at D:E - this pushes the old D:K data to F:M
#    and keeps column A:C (labels) untouched.
$ws.Range("D1:E1").EntireColumn.Insert()

# 2) Copy the cell formatting (number format / alignment / style) from the
#    freshly-shifted F:G columns into the new D:E columns so the new columns
#    look identical to their neighbours (date format on row 7/38/80, right
#    aligned numeric format elsewhere).
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Populate the new D (newest quarter) and E (prior quarter) columns with
#    their values, row by row. $null means the row stays blank (it was blank
#    in the neighbouring columns too), and "NA" writes the literal text used
#    elsewhere in that row.
$newData = @{
    7 = @(43465, 43373)
    8 = @(36900, 36300)
    9 = @(19500, 17900)
    10 = @(17400, 18400)
    11 = @($null, $null)
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(0, 0)
    15 = @(0, 0)
    16 = @($null, $null)
    17 = @(26600, 24700)
    18 = @(10300, 11600)
    19 = @($null, $null)
    20 = @(-11300, -4200)
    21 = @("NA", "NA")
    22 = @(0, 0)
    23 = @(-1100, 7500)
    24 = @(0, 0)
    25 = @(0, 0)
    26 = @(-1100, 7500)
    27 = @(-2200, 6700)
    28 = @(0, 0)
    29 = @(0, 0)
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(11300, 4200)
    33 = @(-2200, 6700)
    34 = @(0, 0)
    35 = @(-2200, 6700)
    38 = @(43465, 43373)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(44700, 53600)
    42 = @(0, 0)
    43 = @(890300, 793500)
    44 = @(0, 0)
    45 = @(0, 0)
    46 = @(0, 0)
    47 = @(3000600, 2830500)
    48 = @(0, 0)
    49 = @(0, 0)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(400, 400)
    53 = @(0, 0)
    54 = @(3971500, 3714600)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(499700, 440700)
    58 = @(1498800, 1636000)
    59 = @(8900, 8700)
    60 = @(0, 0)
    61 = @(497100, 288700)
    62 = @(0, 0)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(3407700, 3120800)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(0, 0)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(563800, 593800)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(-2200, 6700)
    82 = @($null, $null)
    83 = @(0, 0)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(-73100, -197400)
    90 = @($null, $null)
    91 = @(0, 0)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(0, 0)
    95 = @($null, $null)
    96 = @(-12600, -12700)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(64200, 228900)
    101 = @(0, 0)
    102 = @(-8900, 31500)
}

foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    if ($null -ne $vals[0]) {
        $ws.Cells.Item($r, 4).Value = $vals[0]
    }
    if ($null -ne $vals[1]) {
        $ws.Cells.Item($r, 5).Value = $vals[1]
    }
}
